$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "trainingaudio/07_pitapi2.wav"
$ws.Range("B2").Value = "pngimages/07_suitcase.png"

$ws.Range("A3").Value = "trainingaudio/09_tipata2.wav"
$ws.Range("B3").Value = "pngimages/09_plane.png"

$ws.Range("A4").Value = "trainingaudio/27_pakapa1.wav"
$ws.Range("B4").Value = "pngimages/27_kiwi.png"

$ws.Range("A5").Value = "trainingaudio/02_pitito3.wav"
$ws.Range("B5").Value = "pngimages/02_pallet.png"

$ws.Range("A6").Value = "trainingaudio/13_kopopi1.wav"
$ws.Range("B6").Value = "pngimages/13_toast.png"

$ws.Range("A7").Value = "trainingaudio/24_takopa1.wav"
$ws.Range("B7").Value = "pngimages/24_banana.png"
